$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20, column A currently holds the phone number as text ("71277620").
# Convert it to a genuine number to match the new redemption row's format.
$ws.Range("A20").Value = 71277620

# Add a new redemption row (row 21) for phone 71277620 redeeming 100 points.
# Force column A to remain text (matching how this particular row was written
# originally) rather than letting Excel auto-convert the numeric-looking
# string into a number. Clear the formatting afterward so no stray
# number-format style lingers on the cell.
$ws.Range("A21").NumberFormat = "@"
$ws.Range("A21").Value = "71277620"
$ws.Range("A21").ClearFormats()
$ws.Range("B21").Value = 100
$ws.Range("C21").Value = "2025-08-18T17:29:26"
